$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.106.60'
$ws.Range("E2").Value = '  +2.69%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.773.78'
$ws.Range("E3").Value = '  -0.95%  '

# Row 4
$ws.Range("E4").Value = '  +0.22%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '338.42'
$ws.Range("E5").Value = '  -0.54%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.30%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3819'
$ws.Range("E7").Value = '  -3.04%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3416'
$ws.Range("E8").Value = '  -1.35%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.04'
$ws.Range("E9").Value = '  -2.54%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.144'
$ws.Range("E10").Value = '  -4.56%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07369'
$ws.Range("E11").Value = '  -1.75%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.23'
$ws.Range("E12").Value = '  +6.47%  '

# Row 13
$ws.Range("E13").Value = '  +0.29%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.376'
$ws.Range("E14").Value = '  -2.04%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.376'
$ws.Range("E15").Value = '  +3.25%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.777.87'
$ws.Range("E16").Value = '  -0.55%  '

# Row 17
$ws.Range("E17").Value = '  -1.88%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06649'
$ws.Range("E18").Value = '  -0.59%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '82.46'
$ws.Range("E19").Value = '  -2.69%  '

# Row 20
$ws.Range("E20").Value = '  +0.30%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.34'
$ws.Range("E21").Value = '  -2.21%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.390'
$ws.Range("E22").Value = '  -2.16%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.126.85'
$ws.Range("E23").Value = '  +2.77%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.06'
$ws.Range("E24").Value = '  -3.10%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.388'
$ws.Range("E25").Value = '  -0.83%  '

# Row 26
$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.433'
$ws.Range("E26").Value = '  -1.87%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.68'
$ws.Range("E27").Value = '  -2.57%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.400'
$ws.Range("E28").Value = '  -4.22%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '154.14'
$ws.Range("E29").Value = '  -2.51%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.979.00'
$ws.Range("E30").Value = '  -0.55%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '134.38'
$ws.Range("E31").Value = '  -1.26%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.042'
$ws.Range("E32").Value = '  +0.31%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.026'
$ws.Range("E33").Value = '  +0.22%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08912'
$ws.Range("E34").Value = '  +0.82%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.67'
$ws.Range("E35").Value = '  -2.74%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02403'
$ws.Range("E36").Value = '  -0.71%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6825'
$ws.Range("E37").Value = '  +0.02%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.303'
$ws.Range("E38").Value = '  -2.20%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06350'
$ws.Range("E39").Value = '  -2.68%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2156'
$ws.Range("E40").Value = '  -2.77%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.237'
$ws.Range("E41").Value = '  -1.22%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.495'
$ws.Range("E42").Value = '  -8.00%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.190'
$ws.Range("E43").Value = '  -2.29%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.26'
$ws.Range("E44").Value = '  -1.87%  '

# Row 45
$ws.Range("E45").Value = '  +0.26%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6253'
$ws.Range("E46").Value = '  -2.19%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.864'
$ws.Range("E47").Value = '  -0.26%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.46'
$ws.Range("E48").Value = '  +0.06%  '

# Row 49
$ws.Range("E49").Value = '  -3.35%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07498'
$ws.Range("E50").Value = '  +4.71%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.203'
$ws.Range("E51").Value = '  +2.34%  '
